$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "100_2" to "100_1"
$ws.Name = "100_1"

# Remove the old "Total nominations received" sub-header row (row 37). This
# shifts rows 38-42 up to 37-41, which already carry the correct totals for
# confirmed/unconfirmed/withdrawn/rejected/returned.
$ws.Rows("37:37").Delete()

# Relabel column A throughout the breakdown sections to the
# "<Category>, <Status>" style labels, and rename the category headers.
$ws.Range("A6").Value  = "Civilian"
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Confirmed"
$ws.Range("A9").Value  = "     Civilian, Unconfirmed"
$ws.Range("A10").Value = "     Civilian, Withdrawn"
$ws.Range("A11").Value = "     Civilian, Rejected  "
$ws.Range("A12").Value = "     Civilian, Returned"

$ws.Range("A13").Value = "Civilian (lists)"
$ws.Range("A14").Value = "     Civilian (lists), New nominations"
$ws.Range("A15").Value = "     Civilian (lists), Confirmed"
$ws.Range("A16").Value = "     Civilian (lists), Unconfirmed"

$ws.Range("A17").Value = "Air Force"
$ws.Range("A18").Value = "     Air Force, New nominations"
$ws.Range("A19").Value = "     Air Force, Confirmed "
$ws.Range("A20").Value = "     Air Force, Unconfirmed"
$ws.Range("A21").Value = "     Air Force, Returned  "

$ws.Range("A22").Value = "Army"
$ws.Range("A23").Value = "     Army, New nominations"
$ws.Range("A24").Value = "     Army, Confirmed "
$ws.Range("A25").Value = "     Army, Unconfirmed"
$ws.Range("A26").Value = "     Army, Returned  "

$ws.Range("A27").Value = "Navy"
$ws.Range("A28").Value = "     Navy, New nominations"
$ws.Range("A29").Value = "     Navy, Confirmed"
$ws.Range("A30").Value = "     Navy, Unconfirmed"

$ws.Range("A31").Value = "Marine Corps"
$ws.Range("A32").Value = "     Marine Corps, New nominations"
$ws.Range("A33").Value = "     Marine Corps, Confirmed"
$ws.Range("A34").Value = "     Marine Corps, Unconfirmed"
$ws.Range("A35").Value = "     Marine Corps, Returned"

# Row 36 used to be the blank "Summary" header; it becomes the new
# "Total new nominations" row and picks up the value that used to live on
# the (now-deleted) row 37, formatted with the same thousands-separator
# number style as the other Total rows.
$ws.Range("A36").Value = "Total new nominations"
$ws.Range("B36").Value = 51929
$ws.Range("B36").NumberFormat = $ws.Range("B37").NumberFormat
$ws.Range("B36").HorizontalAlignment = $ws.Range("B37").HorizontalAlignment

$ws.Range("A41").Value = "Total returned             "
